# Data Driven Framework course completed except Jenkins
#
# 1) Rename the single existing sheet, add two more sheets (openAccountTest,
#    test_suite) after it, keeping the original sheet active/selected.
# 2) Extend "addCustomerTest" with a new "runMode" column and three more
#    data rows.
# 3) Populate the two new sheets with their header + data rows.

$wb = $excel.ActiveWorkbook

# --- Sheets: rename existing, add the two new ones in order -----------------
$wsAddCustomer = $wb.Worksheets.Item(1)
$wsAddCustomer.Name = "addCustomerTest"

$wsOpenAccount = $wb.Worksheets.Add($null, $wsAddCustomer)
$wsOpenAccount.Name = "openAccountTest"

$wsTestSuite = $wb.Worksheets.Add($null, $wsOpenAccount)
$wsTestSuite.Name = "test_suite"

function Fill-Row($ws, $rowNum, $values, $bold) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $cell = $ws.Cells.Item($rowNum, $i + 1)
        $cell.Value = $values[$i]
        if ($bold) {
            $cell.Font.Bold = $true
        }
    }
}

# --- addCustomerTest: add column E (runMode) + rows 3-5 ---------------------
Fill-Row $wsAddCustomer 1 @("firstName", "lastName", "postCode", "alertText", "runMode") $true
Fill-Row $wsAddCustomer 2 @("Sohaib", "Majeed", "123wp", "Customer added successfully", "Y") $false
Fill-Row $wsAddCustomer 3 @("Daniyal", "Ahmed", "4567xy", "Customer added successfully", "N") $false
Fill-Row $wsAddCustomer 4 @("Kashan", "Ali", "76yrt", "Customer added successfully", "Y") $false
Fill-Row $wsAddCustomer 5 @("Usman", "Shabeer", "89rt", "Customer added successfully", "Y") $false

$wsAddCustomer.Range("E3").Select()

# --- openAccountTest ---------------------------------------------------------
Fill-Row $wsOpenAccount 1 @("customer", "currency") $true
Fill-Row $wsOpenAccount 2 @("Sohaib Majeed", "Rupee") $false

$wsOpenAccount.Range("A2").Select()

# --- test_suite ---------------------------------------------------------------
Fill-Row $wsTestSuite 1 @("TCID", "Runmode") $true
Fill-Row $wsTestSuite 2 @("BankManagerLoginTest", "Y") $false
Fill-Row $wsTestSuite 3 @("AddCustomerTest", "Y") $false
Fill-Row $wsTestSuite 4 @("OpenAccountTest", "N") $false

$wsTestSuite.Range("B4").Select()

# --- restore the original sheet as the active / selected tab ----------------
$wsAddCustomer.Activate()
$wsAddCustomer.Range("E3").Select()
